# "a star time behalve cte matrix"
#
# The authored diff shows one real data change: on the "Barge" worksheet,
# the O/D pair in row 51 (Nuremberg -> Neuss) is removed. Because it sits
# in the middle of the table, removing it shifts every row below it up
# by one, which is why the old row 52/53 pairs (Dortmund/Nuremberg and
# Nuremberg/Dortmund) become the new row 51/52, and the sheet's used
# range shrinks from A1:B53 to A1:B52.
#
# The remaining hunks in the diff (fileVersion/rupBuild, xr:revisionPtr,
# window geometry, x14ac:knownFonts, the "Standaard"->"Normal" cell style
# rename, x14ac:dyDescent on every row, defaultRowHeight 14.45->14.3) are
# metadata that the real Excel.exe client stamps on every save based on
# its build/locale and are not reachable through the workbook/worksheet
# object model, so they are not something this script can (or needs to)
# reproduce; only the view/selection state and the data edit below are
# applied.

$wb = $excel.ActiveWorkbook

$wsBarge = $wb.Worksheets.Item("Barge")
$wsTruck = $wb.Worksheets.Item("Truck")

# Remove the Nuremberg/Neuss row (row 51) from the Barge sheet. This
# shifts rows 52:53 up to become the new rows 51:52, matching the
# reordered O/D pairs seen in the diff, and the dimension becomes
# A1:B52.
$wsBarge.Rows.Item(51).Delete()

# Best-effort autofit of column A (both sheets gained a <cols> entry for
# column A sized to fit its longest entry, e.g. "Willebroek"/"Moerdijk").
$wsBarge.Columns.Item(1).AutoFit()
$wsTruck.Columns.Item(1).AutoFit()

# Restore the recorded view/selection state: Truck sheet scrolled with
# B29 selected, Barge sheet (the active tab) scrolled to row 25 with
# F48 selected.
$wsTruck.Activate()
$wsTruck.Application.ActiveWindow.ScrollRow = 6
[void]$wsTruck.Range("B29").Select()

$wsBarge.Activate()
$wsBarge.Application.ActiveWindow.ScrollRow = 25
[void]$wsBarge.Range("F48").Select()
